$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Tasks Left" figure for the last day (Nov 8) from 7 to 2
$ws.Range("C6").Value = 2

# Update the chart title date from "November 4" to "November 8"
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$chart.ChartTitle.Text = "Burndown Chart  for Sprint 2 as of November 8"

# Update the saved cell selection on the sheet
$ws.Range("K5").Select() | Out-Null
